# Generate Report for Handoff
# - Status text moves from "Handed back: in sync with en-US" to "Ready for handoff"
# - Timestamps are refreshed to the new handoff-generation run
# - The now-shorter "Status" columns are narrowed to fit the new text

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status cells: "Handed back: in sync with en-US" -> "Ready for handoff" ---
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value = "Ready for handoff"
$dede.Range("C2").Value = "Ready for handoff"

# --- Refreshed datetime stamps ---
$overview.Range("G2").Value = "2016-08-18 00:55:42"
$dede.Range("H2").Value = "2016-08-18 00:55:42"
$zhcn.Range("H2").Value = "2016-08-18 00:55:37"

# --- Narrow the Status columns to fit the shorter text ---
$overview.Range("E1").ColumnWidth = 16.333333333333336
$overview.Range("F1").ColumnWidth = 16.333333333333336
$zhcn.Range("C1").ColumnWidth = 16.333333333333336
$dede.Range("C1").ColumnWidth = 16.333333333333336
